$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 14; existing rows 14-19 shift down to 15-20
$ws.Rows.Item(14).Insert()

# Fill in the two new data rows first
$ws.Range("A15").Value = "Report_AzarGH_FeederPeak"
$ws.Range("B15").Value = "2021 May 24"

$ws.Range("A16").Value = "Report_AzarGH_PostFeederLoad"
$ws.Range("B16").Value = "2021 May 24"

# "1400/03/03" looks like a date to the recalculation engine, so typing it
# straight into C15/C16 would silently turn it into a date serial number.
# Enter it with a leading apostrophe in a scratch cell (forcing text, like a
# user pressing ' before typing), then copy just that value into the two
# target cells so their existing cell formatting (borders/fill) is kept.
$ws.Range("Z1").Value = "'1400/03/03"
$ws.Range("Z1").Copy()
$ws.Range("C15").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("C16").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# New header row 14: merge A14:C14 and copy the look of the existing
# "1400 / 02" header row (row 8) onto it
$ws.Range("A14:C14").Merge()
$ws.Range("A8:C8").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A14").Value = "1400 / 03"

# Match the new selection shown in the saved workbook
$ws.Range("A15").Select()
